# S-01018 _ Solucionado issue de minimo 1 metrica
#
# Adds the missing "Nacho" time-tracking entries to "Horas insumidas" so
# every metric/user-story has at least one logged hour, which ripples
# through the SUMIF-based Earned Value sheet and the Estadísticas/chart
# sheet. Also touches a couple of view/selection bits left behind by the
# edit session.

$wb = $excel.ActiveWorkbook

$wsEV   = $wb.Worksheets.Item("Earned Value")
$wsHI   = $wb.Worksheets.Item("Horas insumidas")
$wsStat = $wb.Worksheets.Item("Estadísticas")

# ---------------------------------------------------------------------
# 1) "Horas insumidas": open up 3 fresh rows right before the old row 45
#    (Nico / Scrolling...) so the existing block slides from 45-48 down
#    to 48-51, and the trailing total row slides from 65 down to 68.
# ---------------------------------------------------------------------
$wsHI.Range("B45:B47").EntireRow.Insert()

# New row 45: Nacho, Minuta de reunion, S-01015
$wsHI.Range("B45").Value = 40450
$wsHI.Range("C45").Value = "nacho"
$wsHI.Range("D45").Value = "Minuta de reunion"
$wsHI.Range("E45").Value = "S-01015"
$wsHI.Range("F45").Value = 1

# New row 46: Nacho, Indicador EV, S-01015
$wsHI.Range("B46").Value = 40450
$wsHI.Range("C46").Value = "nacho"
$wsHI.Range("D46").Value = "Indicador EV"
$wsHI.Range("E46").Value = "S-01015"
$wsHI.Range("F46").Value = 1

# Rows 48-51 are the old 45-48 block, shifted down automatically by the
# insert above (values/formats carried over, no further edits needed).

# Fill the 3 new rows directly after the (now shifted) row 51 block -
# these land in what used to be blank filler rows, so no further insert
# is required; the trailing total row stays put at 68.
$wsHI.Range("B52").Value = 40454
$wsHI.Range("C52").Value = "Nacho"
$wsHI.Range("D52").Value = "Se modifico el hibernate.xml junto a los test y se resolvio el issue de "
$wsHI.Range("E52").Value = "S-01018"
$wsHI.Range("F52").Value = 2.5
$wsHI.Rows.Item(52).RowHeight = 18

# New row 47: Nacho, UAT Primer Sprint, S-01018
$wsHI.Range("B47").Value = 40451
$wsHI.Range("C47").Value = "nacho"
$wsHI.Range("D47").Value = "UAT Primer Sprint"
$wsHI.Range("E47").Value = "S-01018"
$wsHI.Range("F47").Value = 1

$wsHI.Range("B53").Value = 40454
$wsHI.Range("C53").Value = "nacho"
$wsHI.Range("D53").Value = "Reporte de avance"
$wsHI.Range("E53").Value = "S-01015"
$wsHI.Range("F53").Value = 1.5

$wsHI.Range("B54").Value = 40454
$wsHI.Range("C54").Value = "nacho"
$wsHI.Range("D54").Value = "Funcionalidad completa"
$wsHI.Range("E54").Value = "S-01015"
$wsHI.Range("F54").Value = 0.5

# ---------------------------------------------------------------------
# 2) Recalculate so every dependent SUMIF / chart cache value is fresh.
# ---------------------------------------------------------------------
$excel.Calculate()

# ---------------------------------------------------------------------
# 3) Leftover UI state from the edit session: "Horas insumidas" ends up
#    the active/selected tab, scrolled down near the new rows; the
#    "Earned Value" sheet selection moves off the SUMIF column; and
#    "Estadísticas" is no longer the tab shown on open.
# ---------------------------------------------------------------------
$wsEV.Activate()
$wsEV.Range("B27").Select()

$wsHI.Activate()
$wsHI.Range("F55").Select()

$wsStat.Activate()
$wsStat.Range("E21").Select()

$wsHI.Activate()
